# Update odds values in row 2 of Sheet1 to reflect the latest FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 3.15
$ws.Range("I2").Value = 3.3
$ws.Range("J2").Value = 2.72
$ws.Range("K2").Value = 2.07
$ws.Range("L2").Value = 3.85
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 7.5
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.39
$ws.Range("T2").Value = 2.8
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 2.25
$ws.Range("W2").Value = 9.25
$ws.Range("X2").Value = 12.5
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 7.5
$ws.Range("AD2").Value = 6.2
$ws.Range("AE2").Value = 11.5
$ws.Range("AI2").Value = 18.5
$ws.Range("AL2").Value = 28
$ws.Range("AO2").Value = 11.25
$ws.Range("AR2").Value = 65
$ws.Range("AT2").Value = 2.8
$ws.Range("AU2").Value = 6.4
$ws.Range("AW2").Value = 5.4
$ws.Range("AX2").Value = 19
